$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 150, shifting existing rows 150:171 down to 151:172.
$ws.Rows.Item(150).Insert()

# Populate the new row 150 with the weekly price-report entry.
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 44491
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = 100112009
$ws.Cells.Item(150, 7).Value = "Acelga"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 100
$ws.Cells.Item(150, 11).Value = 350
$ws.Cells.Item(150, 12).Value = 400
$ws.Cells.Item(150, 13).Value = 375
$ws.Cells.Item(150, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(150, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(150, 16).Value = 375
$ws.Cells.Item(150, 17).Value = 1
$ws.Cells.Item(150, 18).Value = "Hortaliza"
